$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3343.6667
$ws.Range("I11").Value = 3343.6667
$ws.Range("K11").Value = 3343.6667
$ws.Range("M11").Value = -3203.6667
$ws.Range("H15").Value = 1093.04
$ws.Range("I15").Value = 1093.04
$ws.Range("K15").Value = 3279.12
$ws.Range("M15").Value = -3110.12
$ws.Range("H17").Value = 3987.1875
$ws.Range("I17").Value = 219.66667
$ws.Range("J17").Value = 4376.931
$ws.Range("K17").Value = 659.00001
$ws.Range("L17").Value = 13130.793
$ws.Range("M17").Value = -491.00001
$ws.Range("N17").Value = -13466.793
$ws.Range("H52").Value = 400
$ws.Range("J52").Value = 300
$ws.Range("L52").Value = 900
$ws.Range("N52").Value = -1220
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -30972
$ws.Range("H70").Value = 1669.4615
$ws.Range("I70").Value = 817
$ws.Range("J70").Value = 2048.3333
$ws.Range("K70").Value = 2451
$ws.Range("L70").Value = 6144.999899999999
$ws.Range("M70").Value = -2181
$ws.Range("N70").Value = -6684.999899999999
$ws.Range("H73").Value = 1669.4615
$ws.Range("I73").Value = 817
$ws.Range("J73").Value = 2048.3333
$ws.Range("K73").Value = 2451
$ws.Range("L73").Value = 6144.999899999999
$ws.Range("M73").Value = -1515
$ws.Range("N73").Value = -8016.999899999999
$ws.Range("H116").Value = 5446.154
$ws.Range("I116").Value = 4967.143
$ws.Range("K116").Value = 4967.143
$ws.Range("M116").Value = -1525.143
$ws.Range("H137").Value = 92197.89999999999
$ws.Range("I137").Value = 164286.1
$ws.Range("J137").Value = 4090.111
$ws.Range("K137").Value = 492858.3
$ws.Range("L137").Value = 12270.333
$ws.Range("M137").Value = -490308.3
$ws.Range("N137").Value = -17370.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 943.4545000000001
$ws.Range("I2").Value = 723.3158
$ws.Range("K2").Value = 723.3158
$ws.Range("M2").Value = -610.3158
$ws.Range("H5").Value = 794.3043
$ws.Range("I5").Value = 838.82355
$ws.Range("J5").Value = 668.1667
$ws.Range("K5").Value = 838.82355
$ws.Range("L5").Value = 668.1667
$ws.Range("M5").Value = -726.82355
$ws.Range("N5").Value = -892.1667
$ws.Range("H16").Value = 1698
$ws.Range("I16").Value = 1029.3334
$ws.Range("K16").Value = 1029.3334
$ws.Range("M16").Value = -742.3334
$ws.Range("H25").Value = 1999.6666
$ws.Range("I25").Value = 1999.6666
$ws.Range("K25").Value = 1999.6666
$ws.Range("M25").Value = -1597.6666
$ws.Range("H32").Value = 10261.25
$ws.Range("I32").Value = 7097.2856
$ws.Range("K32").Value = 7097.2856
$ws.Range("M32").Value = -6810.2856
$ws.Range("H76").Value = 93057.60000000001
$ws.Range("I76").Value = 149999
$ws.Range("J76").Value = 78822.25
$ws.Range("K76").Value = 149999
$ws.Range("L76").Value = 78822.25
$ws.Range("M76").Value = -149661
$ws.Range("N76").Value = -79498.25
$ws.Range("H79").Value = 93057.60000000001
$ws.Range("I79").Value = 149999
$ws.Range("J79").Value = 78822.25
$ws.Range("K79").Value = 149999
$ws.Range("L79").Value = 78822.25
$ws.Range("M79").Value = -148829
$ws.Range("N79").Value = -81162.25
$ws.Range("H116").Value = 943.4545000000001
$ws.Range("I116").Value = 723.3158
$ws.Range("K116").Value = 723.3158
$ws.Range("M116").Value = 1570.6842
$ws.Range("H132").Value = 2076.2
$ws.Range("I132").Value = 1332.4375
$ws.Range("J132").Value = 2926.2144
$ws.Range("K132").Value = 3997.3125
$ws.Range("L132").Value = 8778.643199999999
$ws.Range("M132").Value = -1467.3125
$ws.Range("N132").Value = -13838.6432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 943.4545000000001
$ws.Range("I3").Value = 723.3158
$ws.Range("K3").Value = 723.3158
$ws.Range("M3").Value = -609.3158
$ws.Range("H4").Value = 794.3043
$ws.Range("I4").Value = 838.82355
$ws.Range("J4").Value = 668.1667
$ws.Range("K4").Value = 838.82355
$ws.Range("L4").Value = 668.1667
$ws.Range("M4").Value = -723.82355
$ws.Range("N4").Value = -898.1667
$ws.Range("H12").Value = 99.75
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 99.5
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 99.5
$ws.Range("M12").Value = 68
$ws.Range("N12").Value = -435.5
$ws.Range("H22").Value = 1716.5834
$ws.Range("I22").Value = 1509
$ws.Range("K22").Value = 1509
$ws.Range("M22").Value = -1336
$ws.Range("H24").Value = 1255
$ws.Range("I24").Value = 1445
$ws.Range("K24").Value = 1445
$ws.Range("M24").Value = -1210
$ws.Range("H29").Value = 217800
$ws.Range("I29").Value = 295000
$ws.Range("J29").Value = 102000
$ws.Range("K29").Value = 295000
$ws.Range("L29").Value = 102000
$ws.Range("M29").Value = -294711
$ws.Range("N29").Value = -102578
$ws.Range("H105").Value = 1711.8
$ws.Range("I105").Value = 1774.3077
$ws.Range("K105").Value = 1774.3077
$ws.Range("M105").Value = -27.30770000000007

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 260.8125
$ws.Range("I7").Value = 30.333334
$ws.Range("J7").Value = 557.1429000000001
$ws.Range("K7").Value = 30.333334
$ws.Range("L7").Value = 557.1429000000001
$ws.Range("M7").Value = 82.66666599999999
$ws.Range("N7").Value = -783.1429000000001
$ws.Range("H17").Value = 2899.4
$ws.Range("J17").Value = 6500
$ws.Range("L17").Value = 6500
$ws.Range("N17").Value = -6848
$ws.Range("H22").Value = 568.8182
$ws.Range("I22").Value = 538.5714
$ws.Range("J22").Value = 621.75
$ws.Range("K22").Value = 538.5714
$ws.Range("L22").Value = 621.75
$ws.Range("M22").Value = -188.5714
$ws.Range("N22").Value = -1321.75
$ws.Range("H31").Value = 13653
$ws.Range("I31").Value = 1379.2307
$ws.Range("K31").Value = 1379.2307
$ws.Range("M31").Value = -1084.2307
$ws.Range("H34").Value = 13653
$ws.Range("I34").Value = 1379.2307
$ws.Range("K34").Value = 1379.2307
$ws.Range("M34").Value = -1177.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3622.8235
$ws.Range("I126").Value = 3324.5
$ws.Range("J126").Value = 4338.8
$ws.Range("K126").Value = 9973.5
$ws.Range("L126").Value = 13016.4
$ws.Range("M126").Value = -7503.5
$ws.Range("N126").Value = -17956.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3404.0454
$ws.Range("I7").Value = 1822.8
$ws.Range("K7").Value = 1822.8
$ws.Range("M7").Value = -1710.8
$ws.Range("H122").Value = 6191.706
$ws.Range("J122").Value = 8221.375
$ws.Range("L122").Value = 24664.125
$ws.Range("N122").Value = -29564.125
$ws.Range("H126").Value = 3404.0454
$ws.Range("I126").Value = 1822.8
$ws.Range("K126").Value = 5468.4
$ws.Range("M126").Value = -2998.4
$ws.Range("H132").Value = 6241.0557
$ws.Range("I132").Value = 6934.5
$ws.Range("K132").Value = 20803.5
$ws.Range("M132").Value = -18273.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9915.333000000001
$ws.Range("J62").Value = 9915.333000000001
$ws.Range("L62").Value = 9915.333000000001
$ws.Range("N62").Value = -11163.333
$ws.Range("H65").Value = 9915.333000000001
$ws.Range("J65").Value = 9915.333000000001
$ws.Range("L65").Value = 49576.665
$ws.Range("N65").Value = -55816.665
$ws.Range("H126").Value = 3064.3928
$ws.Range("I126").Value = 2961.261
$ws.Range("K126").Value = 8883.782999999999
$ws.Range("M126").Value = -6413.782999999999
$ws.Range("H136").Value = 2438.6775
$ws.Range("J136").Value = 4142.4287
$ws.Range("L136").Value = 12427.2861
$ws.Range("N136").Value = -17527.2861
